$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 60, shifting existing rows 60-63 down to 61-64
$ws.Rows.Item(60).Insert()

# Fill in the new row 60 with the new data
$ws.Cells.Item(60, 1).Value = 2
$ws.Cells.Item(60, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(60, 3).Value = "Coquimbo"
$ws.Cells.Item(60, 4).Value = 44714
$ws.Cells.Item(60, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(60, 5).Value = 4
$ws.Cells.Item(60, 6).Value = 100112032
$ws.Cells.Item(60, 7).Value = "Zapallo italiano"
$ws.Cells.Item(60, 8).Value = "Sin especificar"
$ws.Cells.Item(60, 9).Value = "Primera"
$ws.Cells.Item(60, 10).Value = 400
$ws.Cells.Item(60, 11).Value = 10000
$ws.Cells.Item(60, 12).Value = 11000
$ws.Cells.Item(60, 13).Value = 10500
$ws.Cells.Item(60, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(60, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(60, 16).Value = 175
$ws.Cells.Item(60, 17).Value = 60
$ws.Cells.Item(60, 18).Value = "Hortaliza"
